# Applies the "Create Client & Corrected Excels" revision:
#   - Loan gets re-disbursed/re-corrected with a reduced fee (50 instead of 100),
#     which ripples through Summary, Repayment schedule and Transactions sheets.
#   - Stale "Accrual" transaction rows are dropped, collapsing Transactions
#     down to a disbursement / repayment / disbursement trio.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B2").Value = 864.71
$wsSummary.Range("E2").Value = 9135.2900000000009
$wsSummary.Range("F2").Value = 821.7
$wsSummary.Range("A4").Value = 50
$wsSummary.Range("B4").Value = 50

[void]$wsSummary.Activate()
[void]$wsSummary.Range("A7").Select()

# ---------------------------------------------------------------------------
# Repayment schedule sheet
# ---------------------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment schedule")

# Paid Date is now populated for the first installment.
$wsRepay.Range("D3").Value = 42036
$wsRepay.Range("D3").NumberFormat = "d-mmm-yy"

# The placeholder "heading" column stays empty but switches to the italic
# placeholder style used elsewhere on the sheet.
$wsRepay.Range("E3").Value = ""
$wsRepay.Range("E3").Font.Italic = $true

$wsRepay.Range("I3").Value = 50
$wsRepay.Range("K3").Value = 937.72
$wsRepay.Range("P3").Value = 0

[void]$wsRepay.Activate()
[void]$wsRepay.Range("D12").Select()

# ---------------------------------------------------------------------------
# Transactions sheet
# ---------------------------------------------------------------------------
$wsTxn = $wb.Worksheets.Item("Transactions")

# Drop the three stale "Accrual" rows (old rows 2, 3 and 6), leaving the
# disbursement / repayment / disbursement rows which shift up into 2, 3, 4.
$wsTxn.Rows("6:6").Delete()
$wsTxn.Rows("2:3").Delete()

# Row 2 - disbursement, recalculated running loan balance.
$wsTxn.Range("A2").Value = 89
$wsTxn.Range("J2").Value = 9135.2900000000009
$wsTxn.Range("J2").NumberFormat = "#,##0.00"

# Row 3 - repayment, reduced fee + recalculated running loan balance.
$wsTxn.Range("A3").Value = 87
$wsTxn.Range("E3").Value = 937.72
$wsTxn.Range("F3").Value = 864.71
$wsTxn.Range("G3").Value = 23.01
$wsTxn.Range("H3").Value = 50
$wsTxn.Range("J3").Value = 4135.29
$wsTxn.Range("J3").NumberFormat = "#,##0.00"

# Row 4 - original disbursement.
$wsTxn.Range("A4").Value = 84

[void]$wsTxn.Activate()
[void]$wsTxn.Range("C5").Select()
